$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs target cluster)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.434937333333333
$ws.Range("H2").Value = 4.304812
$ws.Range("I2").Value = 0.5010808920723563
$ws.Range("J2").Value = 0.5010808920723562
$ws.Range("M2").Value = 0.3045636666666667
$ws.Range("N2").Value = 0.913691
$ws.Range("Q2").Value = 0.4370297756768889
$ws.Range("R2").Value = 3.933267981092
$ws.Range("S2").Value = 0.5010808920723563
$ws.Range("T2").Value = 0.5010808920723562

# Row 3 (MuSCs target cluster)
$ws.Range("G3").Value = 0.9964423333333334
$ws.Range("I3").Value = 0.3479582011609289
$ws.Range("J3").Value = 0.3479582011609288
$ws.Range("M3").Value = 0.3045636666666667
$ws.Range("N3").Value = 0.913691
$ws.Range("Q3").Value = 0.3034801306618889
$ws.Range("R3").Value = 2.731321175957
$ws.Range("S3").Value = 0.3479582011609289
$ws.Range("T3").Value = 0.3479582011609288

# Row 4 (ECs target cluster)
$ws.Range("G4").Value = 0.4323043333333333
$ws.Range("H4").Value = 1.296913
$ws.Range("I4").Value = 0.150960906766715
$ws.Range("J4").Value = 0.1509609067667149
$ws.Range("M4").Value = 0.3045636666666667
$ws.Range("N4").Value = 0.913691
$ws.Range("Q4").Value = 0.1316641928758889
$ws.Range("R4").Value = 1.184977735883
$ws.Range("S4").Value = 0.150960906766715
$ws.Range("T4").Value = 0.1509609067667149
